# Update the table style applied to every data table in the deck
# (LF Energy "Member Benefits" tables) from the old custom style GUID
# {C320EC70-FFC7-4695-BD8A-B35883984C6E} to the new style GUID
# {50C352CE-1B1D-4D2E-B4E8-68B4B13091BA}.

$p = $ppt.ActivePresentation
$newStyleId = "{50C352CE-1B1D-4D2E-B4E8-68B4B13091BA}"

for ($n = 1; $n -le $p.Slides.Count; $n++) {
    $s = $p.Slides.Item($n)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTable -eq -1) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}
